# feat: add 2022-Q1 data
#
# Before: 2 sheets  -> "2021-Q4" (fund-holding snapshot), "总计" (summary)
# After:  3 sheets  -> "2021-Q4" (unchanged), "2022-Q1" (new fund-holding
#                       snapshot, reusing the old "总计" sheet slot),
#                       "总计" (new sheet appended at the end, summary with
#                       both the 2022-Q1 and 2021-Q4 rows)

$wb = $excel.ActiveWorkbook

$ws2021 = $wb.Worksheets.Item(1)      # "2021-Q4" - stays untouched
$wsTotalOld = $wb.Worksheets.Item(2)  # current "总计" sheet - becomes "2022-Q1"

# ---------------------------------------------------------------------
# 1) Turn the old "总计" sheet into the new "2022-Q1" fund-holding sheet
# ---------------------------------------------------------------------
$wsTotalOld.Cells.Clear()
$wsTotalOld.Name = "2022-Q1"

# Header row: copy format from the existing "2021-Q4" header row (style
# index already bold/centered/bordered) so no new style entries are made.
$ws2021.Range("B1:H1").Copy()
$wsTotalOld.Range("B1").PasteSpecial(-4122)

$wsTotalOld.Range("B1").Value = "基金代码"
$wsTotalOld.Range("C1").Value = "基金名称"
$wsTotalOld.Range("D1").Value = "基金规模"
$wsTotalOld.Range("E1").Value = "股票总仓位"
$wsTotalOld.Range("F1").Value = "仓位占比"
$wsTotalOld.Range("G1").Value = "持有市值(亿元)"
$wsTotalOld.Range("H1").Value = "仓位排名"

# A2 mirrors the style used on "2021-Q4"!A2 (bold/centered/bordered number)
$ws2021.Range("A2").Copy()
$wsTotalOld.Range("A2").PasteSpecial(-4122)
$wsTotalOld.Range("A2").Value = 0

$wsTotalOld.Range("C2").Value = "泰达宏利印度机会股票（QDII）"

# Numeric-looking values stored as TEXT (as in the source data) - force
# text storage via NumberFormat "@" before assignment, otherwise Excel
# auto-coerces "006105" -> 6105 / "0.60" -> 0.6, losing leading/trailing
# zeros and the stored text type.
$wsTotalOld.Range("B2").NumberFormat = "@"
$wsTotalOld.Range("B2").Value = "006105"

$wsTotalOld.Range("D2:G2").NumberFormat = "@"
$wsTotalOld.Range("D2").Value = "0.60"
$wsTotalOld.Range("E2").Value = "87.31"
$wsTotalOld.Range("F2").Value = "3.75"
$wsTotalOld.Range("G2").Value = "0.0225"

$wsTotalOld.Range("H2").Value = 6

# ---------------------------------------------------------------------
# 2) Append a brand-new "总计" sheet at the end with the (old) summary
#    layout, now listing both 2022-Q1 and 2021-Q4.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTotalNew = $wb.Worksheets.Add($null, $lastSheet)
$wsTotalNew.Name = "总计"

# Match the sheetPr/outlinePr + page margins used by the other sheets in
# this workbook (new sheets otherwise get the engine's blank-sheet defaults).
$wsTotalNew.Outline.SummaryRow = 1
$wsTotalNew.Outline.SummaryColumn = 1
$wsTotalNew.PageSetup.LeftMargin = 54
$wsTotalNew.PageSetup.RightMargin = 54
$wsTotalNew.PageSetup.TopMargin = 72
$wsTotalNew.PageSetup.BottomMargin = 72
$wsTotalNew.PageSetup.HeaderMargin = 36
$wsTotalNew.PageSetup.FooterMargin = 36

$wsTotalOld.Range("B1:D1").Copy()
$wsTotalNew.Range("B1").PasteSpecial(-4122)

$wsTotalNew.Range("B1").Value = "日期"
$wsTotalNew.Range("C1").Value = "持有数量(只)"
$wsTotalNew.Range("D1").Value = "持有市值(亿元)"

$wsTotalOld.Range("A2").Copy()
$wsTotalNew.Range("A2").PasteSpecial(-4122)
$wsTotalNew.Range("A2").Value = 0
$wsTotalNew.Range("B2").Value = "2022-Q1"
$wsTotalNew.Range("C2").Value = 1
$wsTotalNew.Range("D2").Value = 0.02

$wsTotalOld.Range("A2").Copy()
$wsTotalNew.Range("A3").PasteSpecial(-4122)
$wsTotalNew.Range("A3").Value = 1
$wsTotalNew.Range("B3").Value = "2021-Q4"
$wsTotalNew.Range("C3").Value = 1
$wsTotalNew.Range("D3").Value = 0.03
